# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" column (E16:E38) previously listed periods in
# ascending order (1805 .. 2003). The account-statement periods were
# refreshed: the old periods were removed and the new set of periods was
# re-entered in descending (most-recent-first) order, so the same 23
# periods now read 2003 down to 1805.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @(
    "2003", "2002", "2001",
    "1912", "1911", "1910", "1909", "1908", "1907", "1906", "1905", "1904",
    "1903", "1902", "1901",
    "1812", "1811", "1810", "1809", "1808", "1807", "1806", "1805"
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("E$row").Value = $periods[$i]
}
